# Applies the parts-list update: adds the SI7006-A20-IM1R humidity/temp
# sensor (HM1) as a new line item, bumps the C2/C4/C6 capacitor count to
# include C7, and refreshes the dependent totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New part row: HM1 / SI7006-A20-IM1R humidity+temp sensor ---
# Copy formatting from an existing, fully-styled data row so the new
# row matches the rest of the table.
$ws.Range("A9:M9").Copy()
$ws.Range("A23:M23").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# (values written in the same order the strings were first introduced)
$ws.Range("G23").Value = "SI7006-A20-IM1R"
$ws.Range("C23").Value = "Low Power Temp/Humidity Sensor I2C"
$ws.Range("H23").Value = "https://onecall.farnell.com/silicon-labs/si7006-a20-im1r/humidity-temp-sensor-dfn-6/dp/3105971"
$ws.Hyperlinks.Add($ws.Range("H23"), "https://onecall.farnell.com/silicon-labs/si7006-a20-im1r/humidity-temp-sensor-dfn-6/dp/3105971")
$ws.Range("A23").Value = "HM1"
$ws.Range("B23").Value = 22
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = "ONECALL"
$ws.Range("F23").Value = 3105971
$ws.Range("I23").Value = 2.44
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 4
$ws.Range("L23").Value = 4
$ws.Range("M23").Formula = "=L23*I23"

# --- 2. C2, C4, C6 -> C2, C4, C6, C7 (one extra 100nF cap added) ---
$ws.Range("A10").Value = "C2, C4, C6, C7"
$ws.Range("D10").Value = 4
$ws.Range("L10").Value = 20

# --- 3. Totals row moves from 23 to 24; quantity subtotal (D) dropped ---
$ws.Range("M26").Copy()
$ws.Range("M24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("M24").Formula = "=SUM(M2:M23)"

# --- 4. VAT total now references the new subtotal row ---
$ws.Range("M26").Formula = "=1.2*M24"

# --- 5. Old unused per-board-count scratch column (D27:D31) cleared ---
$ws.Range("D27").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("D29").ClearContents()
$ws.Range("D30").ClearContents()
$ws.Range("D31").ClearContents()

# --- 6. Selection cosmetic change ---
$ws.Range("L11").Select() | Out-Null

$wb.Save() | Out-Null
